# Auto-generated edit script
$wb = $excel.ActiveWorkbook

# --- sheet1 (Worksheets.Item(1)) ---
$ws = $wb.Worksheets.Item(1)

# Row 97 updates
$ws.Range("B97").Value = 0.04
$ws.Range("D97").Value = 1.03
$ws.Range("H97").Value = 0.36
$ws.Range("I97").Value = -0.65
$ws.Range("J97").Value = -0.04
$ws.Range("L97").Value = -7.92
$ws.Range("M97").Value = -0.13
$ws.Range("O97").Value = -0.14
$ws.Range("Q97").Value = -1.27
$ws.Range("V97").Value = 0.05
$ws.Range("Y97").Value = -0.01
$ws.Range("AC97").Value = -0.2
$ws.Range("AE97").Value = 3.87
$ws.Range("AQ97").Value = -0.29
$ws.Range("AS97").Value = -1.56
$ws.Range("AT97").Value = -2.15
$ws.Range("AU97").Value = -0.14
$ws.Range("AV97").Value = 0.24
$ws.Range("AW97").Value = -0.18
$ws.Range("AX97").Value = 0.59
$ws.Range("AY97").Value = -10.02

# Row 102 updates
$ws.Range("B102").Value = 0.1
$ws.Range("D102").Value = -2.2
$ws.Range("H102").Value = 0.17
$ws.Range("I102").Value = 0.3
$ws.Range("L102").Value = -0.72
$ws.Range("M102").Value = 0.02
$ws.Range("O102").Value = 0.43
$ws.Range("P102").Value = -0.09
$ws.Range("Q102").Value = -1.12
$ws.Range("R102").Value = -0.06
$ws.Range("V102").Value = 0.56
$ws.Range("AG102").Value = 0.02
$ws.Range("AH102").Value = -0.03
$ws.Range("AI102").Value = -0.13
$ws.Range("AS102").Value = 9.67
$ws.Range("AT102").Value = -1.99
$ws.Range("AU102").Value = 0.17
$ws.Range("AW102").Value = 0.09
$ws.Range("AX102").Value = 0.85
$ws.Range("AY102").Value = 8.55

# New row 103
$ws.Range("A103").Value = 45444
$ws.Range("B103").Value = -0.06
$ws.Range("C103").Value = 0.06
$ws.Range("D103").Value = -1.51
$ws.Range("E103").Value = 0.01
$ws.Range("F103").Value = 0.01
$ws.Range("G103").Value = -0.69
$ws.Range("H103").Value = 0.33
$ws.Range("I103").Value = 1.16
$ws.Range("J103").Value = 0.01
$ws.Range("K103").Value = -0.03
$ws.Range("L103").Value = -1.14
$ws.Range("M103").Value = -0.03
$ws.Range("N103").Value = -0.01
$ws.Range("O103").Value = -0.04
$ws.Range("P103").Value = -0.03
$ws.Range("Q103").Value = -2.09
$ws.Range("R103").Value = -0.09
$ws.Range("S103").Value = 0
$ws.Range("T103").Value = 0.26
$ws.Range("U103").Value = -0.03
$ws.Range("V103").Value = 0.18
$ws.Range("W103").Value = -0.65
$ws.Range("X103").Value = -0.14
$ws.Range("Y103").Value = -0.02
$ws.Range("Z103").Value = 0
$ws.Range("AA103").Value = -0.03
$ws.Range("AB103").Value = 0
$ws.Range("AC103").Value = 0.17
$ws.Range("AD103").Value = 0.3
$ws.Range("AE103").Value = 0.79
$ws.Range("AF103").Value = 0.02
$ws.Range("AG103").Value = -0.07
$ws.Range("AH103").Value = 0.47
$ws.Range("AI103").Value = -0.1
$ws.Range("AJ103").Value = 0.01
$ws.Range("AK103").Value = 0.01
$ws.Range("AL103").Value = -0.02
$ws.Range("AM103").Value = 0.01
$ws.Range("AN103").Value = 0.01
$ws.Range("AO103").Value = 0.03
$ws.Range("AP103").Value = -0.02
$ws.Range("AQ103").Value = -0.18
$ws.Range("AR103").Value = -0.04
$ws.Range("AS103").Value = 0.02
$ws.Range("AT103").Value = -6.66
$ws.Range("AU103").Value = 0.06
$ws.Range("AV103").Value = -0.07
$ws.Range("AW103").Value = 0.05
$ws.Range("AX103").Value = 0.02
$ws.Range("AY103").Value = -9.78
$ws.Range("A97").Copy()
$ws.Range("A103").PasteSpecial(-4122)

# --- sheet2 (Worksheets.Item(2)) ---
$ws = $wb.Worksheets.Item(2)

# Row 97 updates
$ws.Range("B97").Value = -0.008
$ws.Range("C97").Value = -0.245
$ws.Range("D97").Value = -0.149
$ws.Range("E97").Value = -0.018
$ws.Range("F97").Value = -0.019
$ws.Range("H97").Value = 0.003
$ws.Range("I97").Value = -0.004
$ws.Range("J97").Value = -0.017
$ws.Range("K97").Value = -0.024
$ws.Range("L97").Value = 0.075
$ws.Range("T97").Value = 0.005
$ws.Range("U97").Value = -0.106
$ws.Range("AA97").Value = -0.008
$ws.Range("AD97").Value = -4.366
$ws.Range("AE97").Value = -0.169
$ws.Range("AF97").Value = -0.009
$ws.Range("AG97").Value = 1.534
$ws.Range("AK97").Value = -0.15
$ws.Range("AL97").Value = 0.324
$ws.Range("AN97").Value = -0.009
$ws.Range("AP97").Value = -0.01
$ws.Range("AQ97").Value = -0.001
$ws.Range("AR97").Value = -0.009
$ws.Range("AS97").Value = 0.009
$ws.Range("AT97").Value = 0.01
$ws.Range("AU97").Value = -0.03
$ws.Range("AV97").Value = -0.023
$ws.Range("AW97").Value = 0.004
$ws.Range("AY97").Value = -0.028
$ws.Range("BA97").Value = 0.037
$ws.Range("BB97").Value = -0.288
$ws.Range("BC97").Value = -0.022
$ws.Range("BD97").Value = 1.061
$ws.Range("BE97").Value = -0.015
$ws.Range("BF97").Value = 0.003
$ws.Range("BJ97").Value = -0.863
$ws.Range("BK97").Value = -0.048
$ws.Range("BL97").Value = -0.481

# Row 102 updates
$ws.Range("B102").Value = 0.097
$ws.Range("C102").Value = 0.049
$ws.Range("D102").Value = 0.033
$ws.Range("E102").Value = 0.038
$ws.Range("H102").Value = 0.005
$ws.Range("I102").Value = -0.007
$ws.Range("J102").Value = 0.011
$ws.Range("K102").Value = 0.035
$ws.Range("L102").Value = 0.067
$ws.Range("N102").Value = 0
$ws.Range("R102").Value = 0.006
$ws.Range("T102").Value = 0
$ws.Range("U102").Value = -0.016
$ws.Range("AD102").Value = -1.978
$ws.Range("AJ102").Value = 0.028
$ws.Range("AK102").Value = -0.073
$ws.Range("AL102").Value = -0.145
$ws.Range("AN102").Value = -0.016
$ws.Range("AO102").Value = -0.005
$ws.Range("AP102").Value = -0.003
$ws.Range("AR102").Value = 0.008
$ws.Range("AS102").Value = -0.009
$ws.Range("AT102").Value = 0.006
$ws.Range("AU102").Value = 0.02
$ws.Range("AV102").Value = -0.002
$ws.Range("AW102").Value = 0.014
$ws.Range("AY102").Value = 0.096
$ws.Range("AZ102").Value = 0.148
$ws.Range("BA102").Value = 0.158
$ws.Range("BB102").Value = 0.881
$ws.Range("BC102").Value = -0.003
$ws.Range("BD102").Value = 0.166
$ws.Range("BF102").Value = -0.006
$ws.Range("BH102").Value = -0.02
$ws.Range("BJ102").Value = 0.206
$ws.Range("BK102").Value = 0.081
$ws.Range("BL102").Value = 0.35
$ws.Range("BM102").Value = 0.005
$ws.Range("BN102").Value = 0.005

# New row 103
$ws.Range("A103").Value = 45444
$ws.Range("B103").Value = -0.048
$ws.Range("C103").Value = -0.074
$ws.Range("D103").Value = -0.106
$ws.Range("E103").Value = -0.009
$ws.Range("F103").Value = -0.01
$ws.Range("G103").Value = 0.006
$ws.Range("H103").Value = 0.002
$ws.Range("I103").Value = 0.005
$ws.Range("J103").Value = -0.015
$ws.Range("K103").Value = -0.065
$ws.Range("L103").Value = 0.07
$ws.Range("M103").Value = 0
$ws.Range("N103").Value = -0.001
$ws.Range("O103").Value = -0.001
$ws.Range("P103").Value = 0
$ws.Range("Q103").Value = -0.007
$ws.Range("R103").Value = 0.003
$ws.Range("S103").Value = -0.004
$ws.Range("T103").Value = -0.002
$ws.Range("U103").Value = 0.039
$ws.Range("V103").Value = -0.007
$ws.Range("W103").Value = -0.004
$ws.Range("X103").Value = -0.133
$ws.Range("Y103").Value = 0.01
$ws.Range("Z103").Value = -0.002
$ws.Range("AA103").Value = -0.001
$ws.Range("AB103").Value = 0
$ws.Range("AC103").Value = 0
$ws.Range("AD103").Value = -2.573
$ws.Range("AE103").Value = -0.344
$ws.Range("AF103").Value = -0.001
$ws.Range("AG103").Value = -1.28
$ws.Range("AH103").Value = -0.198
$ws.Range("AI103").Value = 0.005
$ws.Range("AJ103").Value = -0.028
$ws.Range("AK103").Value = -0.057
$ws.Range("AL103").Value = -0.278
$ws.Range("AM103").Value = 0.001
$ws.Range("AN103").Value = -0.007
$ws.Range("AO103").Value = -0.012
$ws.Range("AP103").Value = -0.011
$ws.Range("AQ103").Value = 0
$ws.Range("AR103").Value = -0.016
$ws.Range("AS103").Value = 0.028
$ws.Range("AT103").Value = -0.004
$ws.Range("AU103").Value = 0.001
$ws.Range("AV103").Value = -0.176
$ws.Range("AW103").Value = 0.039
$ws.Range("AX103").Value = -0.002
$ws.Range("AY103").Value = -0.095
$ws.Range("AZ103").Value = -0.05
$ws.Range("BA103").Value = 0.045
$ws.Range("BB103").Value = 0.959
$ws.Range("BC103").Value = -0.042
$ws.Range("BD103").Value = -0.106
$ws.Range("BE103").Value = 0.007
$ws.Range("BF103").Value = -0.003
$ws.Range("BG103").Value = -0.002
$ws.Range("BH103").Value = -0.001
$ws.Range("BI103").Value = 0.003
$ws.Range("BJ103").Value = -0.057
$ws.Range("BK103").Value = -0.1
$ws.Range("BL103").Value = -0.063
$ws.Range("BM103").Value = -0.133
$ws.Range("BN103").Value = -0.133
$ws.Range("A97").Copy()
$ws.Range("A103").PasteSpecial(-4122)

# --- sheet3 (Worksheets.Item(3)) ---
$ws = $wb.Worksheets.Item(3)

# Row 97 updates
$ws.Range("B97").Value = 0.02
$ws.Range("D97").Value = 0.13
$ws.Range("H97").Value = -5.34
$ws.Range("I97").Value = 5.23
$ws.Range("K97").Value = -0.29
$ws.Range("L97").Value = -2.16
$ws.Range("M97").Value = -0.03
$ws.Range("O97").Value = 0.05
$ws.Range("V97").Value = 0
$ws.Range("X97").Value = -0.78
$ws.Range("Y97").Value = 0.01
$ws.Range("AA97").Value = 0
$ws.Range("AC97").Value = 0.03
$ws.Range("AE97").Value = 1.13
$ws.Range("AG97").Value = -0.01
$ws.Range("AH97").Value = 0.03
$ws.Range("AI97").Value = 0.08
$ws.Range("AL97").Value = 0.11
$ws.Range("AS97").Value = 3.66
$ws.Range("AT97").Value = -4.06
$ws.Range("AU97").Value = -2.32
$ws.Range("AW97").Value = -0.32
$ws.Range("AX97").Value = -0.18
$ws.Range("AY97").Value = -5.04

# Row 102 updates
$ws.Range("D102").Value = 0.09
$ws.Range("I102").Value = 1.43
$ws.Range("L102").Value = -2.24
$ws.Range("M102").Value = -0.03
$ws.Range("O102").Value = -0.78
$ws.Range("P102").Value = -0.02
$ws.Range("Q102").Value = -1.08
$ws.Range("R102").Value = 0.02
$ws.Range("S102").Value = -0.02
$ws.Range("V102").Value = 0.14
$ws.Range("X102").Value = -0.02
$ws.Range("AE102").Value = 0.15
$ws.Range("AG102").Value = 0.16
$ws.Range("AH102").Value = -0.19
$ws.Range("AI102").Value = -0.04
$ws.Range("AK102").Value = -0.02
$ws.Range("AS102").Value = 3.47
$ws.Range("AU102").Value = 0.05
$ws.Range("AV102").Value = -0.2
$ws.Range("AW102").Value = 0.1
$ws.Range("AY102").Value = -10.22

# New row 103
$ws.Range("A103").Value = 45444
$ws.Range("B103").Value = -0.17
$ws.Range("C103").Value = 0
$ws.Range("D103").Value = 0.44
$ws.Range("E103").Value = 0
$ws.Range("F103").Value = 0.02
$ws.Range("G103").Value = 0.42
$ws.Range("H103").Value = -3.69
$ws.Range("I103").Value = 0.52
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = -0.06
$ws.Range("L103").Value = -0.53
$ws.Range("M103").Value = -0.14
$ws.Range("N103").Value = 0
$ws.Range("O103").Value = 0.77
$ws.Range("P103").Value = -0.02
$ws.Range("Q103").Value = -2.68
$ws.Range("R103").Value = -0.08
$ws.Range("S103").Value = -0.01
$ws.Range("T103").Value = -0.01
$ws.Range("U103").Value = 0.02
$ws.Range("V103").Value = 0.1
$ws.Range("W103").Value = -0.91
$ws.Range("X103").Value = 0.55
$ws.Range("Y103").Value = 0
$ws.Range("Z103").Value = 0
$ws.Range("AA103").Value = -0.01
$ws.Range("AB103").Value = -0.01
$ws.Range("AC103").Value = -0.22
$ws.Range("AD103").Value = 0.04
$ws.Range("AE103").Value = 1.59
$ws.Range("AF103").Value = -0.03
$ws.Range("AG103").Value = 0.12
$ws.Range("AH103").Value = -0.02
$ws.Range("AI103").Value = 0.04
$ws.Range("AJ103").Value = 0.01
$ws.Range("AK103").Value = -0.02
$ws.Range("AL103").Value = -0.09
$ws.Range("AM103").Value = -0.01
$ws.Range("AN103").Value = 0
$ws.Range("AO103").Value = 0.09
$ws.Range("AP103").Value = 0.02
$ws.Range("AQ103").Value = -0.14
$ws.Range("AR103").Value = 0.89
$ws.Range("AS103").Value = -0.64
$ws.Range("AT103").Value = -1.76
$ws.Range("AU103").Value = 0.86
$ws.Range("AV103").Value = -0.08
$ws.Range("AW103").Value = 0.15
$ws.Range("AX103").Value = -0.63
$ws.Range("AY103").Value = -5.32
$ws.Range("A97").Copy()
$ws.Range("A103").PasteSpecial(-4122)
